# Natmi following Dr Hou advice
# Update ligand/receptor expressing-cell counts and their derived statistics
# for the Dkk2-Kremen1 sheet (rows 2-4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3, 4)

$values = @{
    2 = @{
        E = 3;  G = 8.398553333333334; H = 25.19566
        K = 3;  M = 3.238669333333334; N = 9.716008
        O = 0.1133140796744973; P = 0.1133140796744973
        Q = 27.20013712503112;  R = 244.80123412528
        S = 0.1133140796744973; T = 0.1133140796744973
    }
    3 = @{
        E = 3;  G = 8.398553333333334; H = 25.19566
        K = 3;  M = 17.985402; N = 53.956206
        O = 0.6292705631384401; P = 0.6292705631384401
        Q = 151.05135791844;     R = 1359.46222126596
        S = 0.6292705631384401; T = 0.6292705631384401
    }
    4 = @{
        E = 3;  G = 8.398553333333334; H = 25.19566
        K = 3;  M = 7.357278333333333; N = 22.071835
        O = 0.2574153571870626; P = 0.2574153571870626
        Q = 61.79049447067779;  R = 556.1144502361
        S = 0.2574153571870626; T = 0.2574153571870626
    }
}

foreach ($r in $rows) {
    $rowVals = $values[$r]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$r").Value = $rowVals[$col]
    }
}
